# Reorder the "Periodo Mora" column (E16:E22) into ascending chronological
# order (YYMM), matching the refreshed EC database:
#   before: 2503, 2502, 2501, 2412, 2411, 2410, 2504
#   after : 2410, 2411, 2412, 2501, 2502, 2503, 2504

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2410"
$ws.Range("E17").Value = "2411"
$ws.Range("E18").Value = "2412"
$ws.Range("E19").Value = "2501"
$ws.Range("E20").Value = "2502"
$ws.Range("E21").Value = "2503"
$ws.Range("E22").Value = "2504"
